# Weekly update: insert a new data row for "Terminal Hortofrutícola Agro
# Chillán - Zanahoria" just before the existing row 240, shifting every
# later row down by one (old row 240 -> new row 241, ... old row 340 ->
# new row 341).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 240; Excel shifts rows 240..340 down to
# 241..341 automatically (this also grows dimension to A1:R341).
$ws.Rows(240).Insert()

# Populate the newly inserted row 240 with this week's record.
$ws.Range("A240").Value = 7
$ws.Range("B240").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C240").Value = "Ñuble"
$ws.Range("D240").Value = 44845
$ws.Range("E240").Value = 16
$ws.Range("F240").Value = 100114013
$ws.Range("G240").Value = "Zanahoria"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 120
$ws.Range("K240").Value = 14000
$ws.Range("L240").Value = 15000
$ws.Range("M240").Value = 14500
$ws.Range("N240").Value = "`$/saco 20 kilos"
$ws.Range("O240").Value = "Región de Ñuble"
$ws.Range("P240").Value = 725
$ws.Range("Q240").Value = 20
$ws.Range("R240").Value = "Hortaliza"
